# Apply scheduled-runner updates to Sheets/Twintania_Profits.xlsx
# Values below come from the source-of-truth diff for each (sheet,row).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 49
$ws.Range("H49").Value = 870.25
$ws.Range("I49").Value = 923.5
$ws.Range("K49").Value = 2770.5
$ws.Range("M49").Value = -2634.5
# Row 80
$ws.Range("H80").Value = 401445.16
$ws.Range("I80").Value = 826.0833
$ws.Range("J80").Value = 771247.4
$ws.Range("K80").Value = 2478.2499
$ws.Range("L80").Value = 2313742.2
$ws.Range("M80").Value = -1480.2499
$ws.Range("N80").Value = -2315738.2
# Row 83
$ws.Range("H83").Value = 401445.16
$ws.Range("I83").Value = 826.0833
$ws.Range("J83").Value = 771247.4
$ws.Range("K83").Value = 7434.7497
$ws.Range("L83").Value = 6941226.600000001
$ws.Range("M83").Value = -2442.7497
$ws.Range("N83").Value = -6951210.600000001
# Row 88
$ws.Range("H88").Value = 2300.4285
$ws.Range("I88").Value = 2099.5
$ws.Range("K88").Value = 2099.5
$ws.Range("M88").Value = -1693.5
# Row 91
$ws.Range("H91").Value = 2300.4285
$ws.Range("I91").Value = 2099.5
$ws.Range("K91").Value = 2099.5
$ws.Range("M91").Value = -695.5
# Row 92
$ws.Range("H92").Value = 667.58826
$ws.Range("I92").Value = 712.61536
$ws.Range("K92").Value = 712.61536
$ws.Range("M92").Value = 535.38464
# Row 135
$ws.Range("H135").Value = 1680.9445
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
# Row 137
$ws.Range("H137").Value = 9591.034
$ws.Range("J137").Value = 17877.273
$ws.Range("L137").Value = 53631.819
$ws.Range("N137").Value = -58731.819
# Row 138
$ws.Range("H138").Value = 4170.7812
$ws.Range("I138").Value = 1752.5
$ws.Range("J138").Value = 5270
$ws.Range("K138").Value = 5257.5
$ws.Range("L138").Value = 15810
$ws.Range("M138").Value = -117.5
$ws.Range("N138").Value = -26090
# Row 141
$ws.Range("H141").Value = 2318.724
$ws.Range("I141").Value = 2557.6667
$ws.Range("J141").Value = 1171.8
$ws.Range("K141").Value = 7673.000100000001
$ws.Range("L141").Value = 3515.4
$ws.Range("M141").Value = -2493.000100000001
$ws.Range("N141").Value = -13875.4

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 1779.81
$ws.Range("I32").Value = 1581.4592
$ws.Range("J32").Value = 11499
$ws.Range("K32").Value = 1581.4592
$ws.Range("L32").Value = 11499
$ws.Range("M32").Value = -1294.4592
$ws.Range("N32").Value = -12073
# Row 74
$ws.Range("H74").Value = 7701.0815
$ws.Range("I74").Value = 8387.098
$ws.Range("K74").Value = 8387.098
$ws.Range("M74").Value = -7513.098
# Row 77
$ws.Range("H77").Value = 7701.0815
$ws.Range("I77").Value = 8387.098
$ws.Range("K77").Value = 41935.49
$ws.Range("M77").Value = -37567.49
# Row 80
$ws.Range("H80").Value = 58384.715
$ws.Range("J80").Value = 58384.715
$ws.Range("L80").Value = 58384.715
$ws.Range("N80").Value = -60380.715
# Row 83
$ws.Range("H83").Value = 58384.715
$ws.Range("J83").Value = 58384.715
$ws.Range("L83").Value = 175154.145
$ws.Range("N83").Value = -185138.145
# Row 122
$ws.Range("H122").Value = 2561.3684
$ws.Range("I122").Value = 2525.889
$ws.Range("J122").Value = 3200
$ws.Range("K122").Value = 7577.667
$ws.Range("L122").Value = 9600
$ws.Range("M122").Value = -5127.667
$ws.Range("N122").Value = -14500
# Row 125
$ws.Range("H125").Value = 75000
$ws.Range("J125").Value = 75000
$ws.Range("L125").Value = 75000
$ws.Range("N125").Value = -84840

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 325552.4
$ws.Range("I86").Value = 478747.34
$ws.Range("J86").Value = 3843.1
$ws.Range("K86").Value = 478747.34
$ws.Range("L86").Value = 3843.1
$ws.Range("M86").Value = -477624.34
$ws.Range("N86").Value = -6089.1
# Row 89
$ws.Range("H89").Value = 325552.4
$ws.Range("I89").Value = 478747.34
$ws.Range("J89").Value = 3843.1
$ws.Range("K89").Value = 2393736.7
$ws.Range("L89").Value = 19215.5
$ws.Range("M89").Value = -2388120.7
$ws.Range("N89").Value = -30447.5

$ws = $wb.Worksheets.Item("CRP")
# Row 25
$ws.Range("H25").Value = 4860.143
$ws.Range("J25").Value = 5802
$ws.Range("L25").Value = 5802
$ws.Range("N25").Value = -6150
# Row 31
$ws.Range("H31").Value = 327365.1
$ws.Range("I31").Value = 145414.58
$ws.Range("J31").Value = 374537.44
$ws.Range("K31").Value = 145414.58
$ws.Range("L31").Value = 374537.44
$ws.Range("M31").Value = -145119.58
$ws.Range("N31").Value = -375127.44
# Row 34
$ws.Range("H34").Value = 327365.1
$ws.Range("I34").Value = 145414.58
$ws.Range("J34").Value = 374537.44
$ws.Range("K34").Value = 145414.58
$ws.Range("L34").Value = 374537.44
$ws.Range("M34").Value = -145212.58
$ws.Range("N34").Value = -374941.44
# Row 132
$ws.Range("H132").Value = 25137.125
$ws.Range("I132").Value = 15509.348
$ws.Range("K132").Value = 46528.044
$ws.Range("M132").Value = -43998.044

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 2707.375
$ws.Range("I5").Value = 1959.6
$ws.Range("J5").Value = 3047.2727
$ws.Range("K5").Value = 5878.799999999999
$ws.Range("L5").Value = 9141.8181
$ws.Range("M5").Value = -5766.799999999999
$ws.Range("N5").Value = -9365.8181
# Row 35
$ws.Range("H35").Value = 450
$ws.Range("I35").Value = 100
$ws.Range("J35").Value = 800
$ws.Range("K35").Value = 300
$ws.Range("L35").Value = 2400
$ws.Range("M35").Value = -12
$ws.Range("N35").Value = -2976
# Row 46
$ws.Range("H46").Value = 598.6667
$ws.Range("J46").Value = 698.25
$ws.Range("L46").Value = 2094.75
$ws.Range("N46").Value = -2276.75
# Row 58
$ws.Range("H58").Value = 5850
$ws.Range("I58").Value = 5775
$ws.Range("K58").Value = 17325
$ws.Range("M58").Value = -17197
# Row 132
$ws.Range("H132").Value = 1899.125
$ws.Range("I132").Value = 1899.125
$ws.Range("K132").Value = 17092.125
$ws.Range("M132").Value = -14562.125
# Row 135
$ws.Range("H135").Value = 2707.375
$ws.Range("I135").Value = 1959.6
$ws.Range("J135").Value = 3047.2727
$ws.Range("K135").Value = 17636.4
$ws.Range("L135").Value = 27425.4543
$ws.Range("M135").Value = -15101.4
$ws.Range("N135").Value = -32495.4543

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 3694.2856
$ws.Range("I80").Value = 3490
$ws.Range("J80").Value = 3966.6667
$ws.Range("K80").Value = 3490
$ws.Range("L80").Value = 3966.6667
$ws.Range("M80").Value = -2492
$ws.Range("N80").Value = -5962.6667
# Row 83
$ws.Range("H83").Value = 3694.2856
$ws.Range("I83").Value = 3490
$ws.Range("J83").Value = 3966.6667
$ws.Range("K83").Value = 17450
$ws.Range("L83").Value = 19833.3335
$ws.Range("M83").Value = -12458
$ws.Range("N83").Value = -29817.3335
# Row 135
$ws.Range("H135").Value = 50000
$ws.Range("J135").Value = 50000
$ws.Range("L135").Value = 50000
$ws.Range("N135").Value = -60140

$ws = $wb.Worksheets.Item("LTW")
# Row 43
$ws.Range("H43").Value = 15000.5
$ws.Range("I43").Value = 15000.5
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 15000.5
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -14807.5
$ws.Range("N43").ClearContents()
# Row 46
$ws.Range("H46").Value = 1713.05
$ws.Range("I46").Value = 1048.75
$ws.Range("J46").Value = 1879.125
$ws.Range("K46").Value = 1048.75
$ws.Range("L46").Value = 1879.125
$ws.Range("M46").Value = -860.75
$ws.Range("N46").Value = -2255.125
# Row 51
$ws.Range("H51").Value = 41857.332
$ws.Range("J51").Value = 41857.332
$ws.Range("L51").Value = 41857.332
$ws.Range("N51").Value = -42813.332
# Row 68
$ws.Range("H68").Value = 2876.0293
$ws.Range("I68").Value = 2565.1035
$ws.Range("J68").Value = 4679.4
$ws.Range("K68").Value = 2565.1035
$ws.Range("L68").Value = 4679.4
$ws.Range("M68").Value = -1816.1035
$ws.Range("N68").Value = -6177.4
# Row 71
$ws.Range("H71").Value = 2876.0293
$ws.Range("I71").Value = 2565.1035
$ws.Range("J71").Value = 4679.4
$ws.Range("K71").Value = 12825.5175
$ws.Range("L71").Value = 23397
$ws.Range("M71").Value = -9081.517500000002
$ws.Range("N71").Value = -30885
# Row 132
$ws.Range("H132").Value = 4451.8096
$ws.Range("I132").Value = 3864.182
$ws.Range("K132").Value = 11592.546
$ws.Range("M132").Value = -9062.545999999998
# Row 136
$ws.Range("H136").Value = 2425.5588
$ws.Range("I136").Value = 1994.1177
$ws.Range("K136").Value = 5982.3531
$ws.Range("M136").Value = -3432.3531

$ws = $wb.Worksheets.Item("WVR")
# Row 38
$ws.Range("H38").Value = 10000
$ws.Range("I38").Value = 10000
$ws.Range("K38").Value = 10000
$ws.Range("M38").Value = -9527
# Row 49
$ws.Range("H49").Value = 99999
$ws.Range("I49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("M49").ClearContents()
# Row 62
$ws.Range("H62").Value = 3872.4375
$ws.Range("I62").Value = 3855.2144
$ws.Range("J62").Value = 3993
$ws.Range("K62").Value = 3855.2144
$ws.Range("L62").Value = 3993
$ws.Range("M62").Value = -3231.2144
$ws.Range("N62").Value = -5241
# Row 65
$ws.Range("H65").Value = 3872.4375
$ws.Range("I65").Value = 3855.2144
$ws.Range("J65").Value = 3993
$ws.Range("K65").Value = 19276.072
$ws.Range("L65").Value = 3993
$ws.Range("M65").Value = -16156.072
$ws.Range("N65").Value = -26205
# Row 126
$ws.Range("H126").Value = 17907.182
$ws.Range("I126").Value = 17907.182
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 53721.546
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -51251.546
$ws.Range("N126").ClearContents()
# Row 132
$ws.Range("H132").Value = 26190.334
$ws.Range("I132").Value = 25475.432
$ws.Range("K132").Value = 76426.296
$ws.Range("M132").Value = -73896.296
# Row 138
$ws.Range("H138").Value = 59047.668
$ws.Range("J138").Value = 59047.668
$ws.Range("L138").Value = 59047.668
$ws.Range("N138").Value = -69327.66800000001
